$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "94.193.66"
Set-TextValue $ws "E2" "  +1.06%  "

Set-TextValue $ws "D3" "3.071.11"
Set-TextValue $ws "E3" "  -1.83%  "

Set-TextValue $ws "E4" "  -0.01%  "

Set-TextValue $ws "D5" "233.53"
Set-TextValue $ws "E5" "  -3.98%  "

Set-TextValue $ws "D6" "607.46"
Set-TextValue $ws "E6" "  -1.62%  "

Set-TextValue $ws "E7" "  -1.12%  "

Set-TextValue $ws "E8" "  -7.85%  "

Set-TextValue $ws "E9" "  -0.06%  "

Set-TextValue $ws "D10" "0.802"
Set-TextValue $ws "E10" "  +9.09%  "

Set-TextValue $ws "D11" "3.071.90"
Set-TextValue $ws "E11" "  -1.64%  "

Set-TextValue $ws "D12" "0.195"
Set-TextValue $ws "E12" "  -3.88%  "

Set-TextValue $ws "D13" "94.038.76"
Set-TextValue $ws "E13" "  +1.36%  "

Set-TextValue $ws "D14" "0.0000239"
Set-TextValue $ws "E14" "  -6.77%  "

Set-TextValue $ws "D15" "33.48"
Set-TextValue $ws "E15" "  -3.31%  "

Set-TextValue $ws "D16" "5.25"
Set-TextValue $ws "E16" "  -4.75%  "

Set-TextValue $ws "D17" "3.647.23"
Set-TextValue $ws "E17" "  -1.71%  "

Set-TextValue $ws "D18" "3.082.15"
Set-TextValue $ws "E18" "  -0.60%  "

Set-TextValue $ws "D19" "3.54"
Set-TextValue $ws "E19" "  -6.48%  "

Set-TextValue $ws "D20" "14.41"
Set-TextValue $ws "E20" "  -2.66%  "

Set-TextValue $ws "D21" "5.68"
Set-TextValue $ws "E21" "  -2.29%  "

Set-TextValue $ws "D22" "438.51"
Set-TextValue $ws "E22" "  -2.77%  "

Set-TextValue $ws "D23" "8.76"
Set-TextValue $ws "E23" "  -7.46%  "

Set-TextValue $ws "D24" "0.0000189"
Set-TextValue $ws "E24" "  -8.64%  "

Set-TextValue $ws "E25" "  +4.76%  "

Set-TextValue $ws "D26" "5.46"
Set-TextValue $ws "E26" "  -6.52%  "

Set-TextValue $ws "D27" "84.67"
Set-TextValue $ws "E27" "  -2.86%  "

Set-TextValue $ws "D28" "11.76"
Set-TextValue $ws "E28" "  -0.93%  "

Set-TextValue $ws "D29" "3.251.83"
Set-TextValue $ws "E29" "  -0.98%  "

Set-TextValue $ws "D31" "0.247"
Set-TextValue $ws "E31" "  +7.98%  "

Set-TextValue $ws "D32" "0.175"
Set-TextValue $ws "E32" "  +3.58%  "

Set-TextValue $ws "D33" "0.119"
Set-TextValue $ws "E33" "  -12.97%  "

Set-TextValue $ws "D34" "9.05"
Set-TextValue $ws "E34" "  -2.36%  "

Set-TextValue $ws "D35" "0.998"
Set-TextValue $ws "E35" "  -0.14%  "

Set-TextValue $ws "D36" "7.68"
Set-TextValue $ws "E36" "  -4.52%  "

Set-TextValue $ws "D37" "0.155"
Set-TextValue $ws "E37" "  -5.87%  "

Set-TextValue $ws "D38" "25.40"
Set-TextValue $ws "E38" "  -3.37%  "

Set-TextValue $ws "E39" "  -2.39%  "

Set-TextValue $ws "B40" "WhiteBITCoin"
Set-TextValue $ws "C40" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws "D40" "23.99"
Set-TextValue $ws "E40" "  +4.00%  "

Set-TextValue $ws "B41" "PolygonEcosystemToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws "D41" "0.438"
Set-TextValue $ws "E41" "  +0.14%  "

Set-TextValue $ws "D42" "467.71"
Set-TextValue $ws "E42" "  -3.81%  "

Set-TextValue $ws "D43" "3.72"
Set-TextValue $ws "E43" "  -10.96%  "

Set-TextValue $ws "D44" "1.24"
Set-TextValue $ws "E44" "  -5.41%  "

Set-TextValue $ws "E45" "  -0.01%  "

Set-TextValue $ws "D46" "3.11"
Set-TextValue $ws "E46" "  -11.49%  "

Set-TextValue $ws "D47" "159.94"
Set-TextValue $ws "E47" "  -0.63%  "

Set-TextValue $ws "E48" "  -5.32%  "

Set-TextValue $ws "D49" "0.669"
Set-TextValue $ws "E49" "  -3.86%  "

Set-TextValue $ws "D50" "43.65"
Set-TextValue $ws "E50" "  -1.04%  "

Set-TextValue $ws "D51" "0.997"
Set-TextValue $ws "E51" "  -0.06%  "
